$d = $word.ActiveDocument

# 1. "Out-x and out-y" -> "Out-1 and out-2" (the two outputs used to drive the H-bridge)
$d.Content.Find.Execute("Out-x and out-y are never allowed", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Out-1 and out-2 are never allowed", 2)

# 2. First placeholder "****" -> "push-sensor" (name of the sensor used at the top position)
$d.Content.Find.Execute("we use a ****.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "we use a push-sensor.", 2)

# 3. Second placeholder "*****" -> "push-sensor"
$d.Content.Find.Execute("this ***** is pressed", $true, $false, $false, $false, $false,
                         $true, 1, $false, "this push-sensor is pressed", 2)

# 4. Word keeps an auto "_GoBack" bookmark at the location of the last edit. Move it to
#    right after the last edit made above (just before " is pressed we immediately ...").
$goBack = $d.Content
$goBack.Find.Execute("push-sensor is pressed", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
$markPos = $goBack.Start + "push-sensor".Length
$markRange = $d.Range($markPos, $markPos)
$d.Bookmarks.Add("_GoBack", $markRange)
